# Auto-generated edit script for horarios-141 workbook update
# Commit: Horarios actualizados Linea 141 - 731

$wb = $excel.ActiveWorkbook

# ----- LP1912 -----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 06:52:31"
$ws.Range("A3").Value = "Total filas: 61"

$rows = @(
  @(24, "04:17:03", "05:35", "215B_EL PATO", 78, "LP1912"),
  @(25, "03:42:43", "05:35", "14_ABASTO", 113, "LP1912"),
  @(41, "06:52:31", "06:59", "14_ABASTO", 7, "LP1912"),
  @(42, "06:52:31", "07:01", "16_SANTA ANA", 9, "LP1912"),
  @(44, "06:52:31", "07:05", "23_HERNANDEZ", 13, "LP1912"),
  @(45, "06:52:31", "07:05", "15_ABASTO", 13, "LP1912"),
  @(46, "05:27:50", "07:06", "225_GOMEZ", 99, "LP1912"),
  @(47, "06:52:31", "07:07", "225_GOMEZ", 15, "LP1912"),
  @(48, "06:52:31", "07:11", "215A_EL PATO", 19, "LP1912"),
  @(49, "06:52:31", "07:15", "11_ETCHEVERRY", 23, "LP1912"),
  @(50, "06:52:31", "07:16", "16_SANTA ANA", 24, "LP1912"),
  @(51, "06:52:31", "07:21", "26_HERNANDEZ", 29, "LP1912"),
  @(52, "06:52:31", "07:23", "10_OLMOS", 31, "LP1912"),
  @(53, "05:55:46", "07:31", "16_SANTA ANA", 96, "LP1912"),
  @(54, "06:52:31", "07:31", "11_ETCHEVERRY", 39, "LP1912"),
  @(55, "06:52:31", "07:32", "84_COLONIA URQUIZA-ESC 49", 40, "LP1912"),
  @(56, "06:52:31", "07:36", "27_EL RETIRO", 44, "LP1912"),
  @(57, "06:52:31", "07:39", "10_OLMOS", 47, "LP1912"),
  @(58, "06:52:31", "07:47", "14_ABASTO", 55, "LP1912"),
  @(59, "06:52:31", "07:51", "215D_EL PATO", 59, "LP1912"),
  @(60, "06:52:31", "08:12", "15_ABASTO", 80, "LP1912"),
  @(61, "06:52:31", "08:21", "26_HERNANDEZ", 89, "LP1912"),
  @(62, "06:52:31", "08:22", "16_P MOR-SANTA ANA", 90, "LP1912"),
  @(63, "06:52:31", "08:23", "215B_EL PATO", 91, "LP1912"),
  @(64, "06:52:31", "08:27", "84_COLONIA URQUIZA-ESC 49", 95, "LP1912"),
  @(65, "06:52:31", "08:35", "23_HERNANDEZ", 103, "LP1912"),
  @(66, "06:52:31", "08:42", "81_EL PELIGRO", 110, "LP1912")
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value = $r[1]
  $ws.Cells.Item($rowNum, 2).Value = $r[2]
  $ws.Cells.Item($rowNum, 3).Value = $r[3]
  $ws.Cells.Item($rowNum, 4).Value = $r[4]
  $ws.Cells.Item($rowNum, 5).Value = $r[5]
}

# ----- LP1912-215 -----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 06:52:31"

$rows = @(
  @(20, "06:52:31", "07:11", "215A_EL PATO", 19, "LP1912"),
  @(21, "06:52:31", "07:51", "215D_EL PATO", 59, "LP1912"),
  @(22, "06:52:31", "08:23", "215B_EL PATO", 91, "LP1912")
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value = $r[1]
  $ws.Cells.Item($rowNum, 2).Value = $r[2]
  $ws.Cells.Item($rowNum, 3).Value = $r[3]
  $ws.Cells.Item($rowNum, 4).Value = $r[4]
  $ws.Cells.Item($rowNum, 5).Value = $r[5]
}

# ----- 6203-6173 -----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 06:52:31"
$ws.Range("A3").Value = "Total filas: 14"

$rows = @(
  @(14, "06:52:31", "07:00", "215B_LP-P MOR-1 Y 57", 8, "L6173"),
  @(17, "06:52:31", "07:41", "215A_LA PLATA", 49, "L6173"),
  @(18, "06:52:31", "08:07", "215C_LA PLATA", 75, "L6203"),
  @(19, "06:52:31", "08:30", "215A_LA PLATA", 98, "L6173")
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value = $r[1]
  $ws.Cells.Item($rowNum, 2).Value = $r[2]
  $ws.Cells.Item($rowNum, 3).Value = $r[3]
  $ws.Cells.Item($rowNum, 4).Value = $r[4]
  $ws.Cells.Item($rowNum, 5).Value = $r[5]
}
